$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text-preserving numeric-looking values in column D need NumberFormat = "@"
# so that values like "0.0002000" keep their exact text (trailing zeros,
# no scientific notation) instead of being auto-converted to a Double by the
# Excel input parser (matches the original inlineStr/text storage of col D).

$priceCells = @("D2","D3","D4","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D22","D23","D24","D25","D26","D27","D40","D41","D42","D43","D44","D45","D47","D49","D50")
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# --- Column D (Price) updates ---
$ws.Range("D2").Value = "247.28"
$ws.Range("D3").Value = "26.36"
$ws.Range("D4").Value = "5.091"
$ws.Range("D7").Value = "3.057"
$ws.Range("D8").Value = "0.8125"
$ws.Range("D9").Value = "0.8428"
$ws.Range("D10").Value = "0.1345"
$ws.Range("D11").Value = "0.03164"
$ws.Range("D12").Value = "0.02818"
$ws.Range("D13").Value = "0.09408"
$ws.Range("D14").Value = "0.001512"
$ws.Range("D15").Value = "0.0005989"
$ws.Range("D16").Value = "0.006252"
$ws.Range("D17").Value = "3.578"
$ws.Range("D18").Value = "2.118"
$ws.Range("D19").Value = "0.3181"
$ws.Range("D20").Value = "0.06959"
$ws.Range("D22").Value = "3.767"
$ws.Range("D23").Value = "0.04675"
$ws.Range("D24").Value = "0.1375"
$ws.Range("D25").Value = "0.001249"
$ws.Range("D26").Value = "0.004623"
$ws.Range("D27").Value = "0.00009599"
$ws.Range("D40").Value = "0.03666"
$ws.Range("D41").Value = "0.006109"
$ws.Range("D42").Value = "0.1061"
$ws.Range("D43").Value = "0.002597"
$ws.Range("D44").Value = "0.008690"
$ws.Range("D45").Value = "0.00005290"
$ws.Range("D47").Value = "0.1200"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("D50").Value = "0.0002000"

# --- Columns B (Coin), C (Link), E (Volume) updates: row re-shuffle + price refresh ---
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("E7").Value = "6GateTokenGT"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("E8").Value = "7MXTokenMX"
$ws.Range("B9").Value = "FTXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("E9").Value = "8FTXTokenFTT"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("E12").Value = "11BitrueCoinBTR"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("E13").Value = "12BitMartTokenBMX"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("E14").Value = "13BitForexTokenBF"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("E15").Value = "14OneONEWorstin24h"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("E16").Value = "15TigerCashTCH"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("E17").Value = "16LEOLEO"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("E18").Value = "17BTSETokenBTSE"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("E19").Value = "18BitpandaEcosystemTokenBEST"
$ws.Range("B20").Value = "MandalaExchangeToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("E20").Value = "19MandalaExchangeTokenMDX"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"
